# upload-product.xlsx: add a Vietnamese "Hướng dẫn" (instructions) sheet
# after "Products", plus a new LEVEL column header on Products!I1.

$wb = $excel.ActiveWorkbook
$products = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Add the new sheet right after "Products".
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Add($null, $products)
$ws.Name = "Hướng dẫn"

# Narrow "bullet" column.
$ws.Columns.Item(2).ColumnWidth = 4.25

# ---------------------------------------------------------------------
# Helpers: every text cell in this sheet is Times New Roman, text-formatted
# (numFmt "@"), vertical-centered. Headings ("- XXX: ...") are bold size 12
# black; sub bullets are size 12 not-bold black (optionally indented);
# the two bold "STOCK_QUANTITY"/"LEVEL" headings late in the sheet use the
# plain size-11 bold font instead.
# ---------------------------------------------------------------------
function Set-Heading12($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Font.Name = "Times New Roman"
    $c.Font.Family = 1
    $c.Font.Size = 12
    $c.Font.Bold = $true
    $c.Font.Color = 0
    $c.NumberFormat = "@"
    $c.VerticalAlignment = -4108
}

function Set-Body12($addr, $text, $indent) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Font.Name = "Times New Roman"
    $c.Font.Family = 1
    $c.Font.Size = 12
    $c.Font.Color = 0
    $c.NumberFormat = "@"
    $c.VerticalAlignment = -4108
    if ($indent -gt 0) {
        $c.HorizontalAlignment = -4131
        $c.IndentLevel = $indent
    }
}

function Set-Heading11($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Font.Name = "Times New Roman"
    $c.Font.Family = 1
    $c.Font.Size = 11
    $c.Font.Bold = $true
    $c.NumberFormat = "@"
}

function Set-Blank12($addr) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Times New Roman"
    $c.Font.Family = 1
    $c.Font.Size = 12
    $c.Font.Bold = $true
    $c.Font.Color = 0
    $c.NumberFormat = "@"
    $c.VerticalAlignment = -4108
}

function Set-BlankIndent($addr) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Times New Roman"
    $c.Font.Family = 1
    $c.Font.Size = 12
    $c.Font.Color = 0
    $c.NumberFormat = "@"
    $c.VerticalAlignment = -4108
    $c.HorizontalAlignment = -4131
    $c.IndentLevel = 1
}

# ---------------------------------------------------------------------
# 2. Fill the sheet. Order matches the original authoring order so the
#    shared-string table indices line up: the explanatory paragraph
#    block was typed first (rows 3-7, 12-17, 18-19, 21-26), then rows
#    2, 8, 9, 10, 11 were inserted afterwards.
# ---------------------------------------------------------------------
Set-Body12    "B3"  "- SHORT_DESC: Mô tả ngắn (Có thể điền hoặc không điền). " 0
Set-Body12    "B4"  "- DESC: Mô tả (Có thể điền hoặc không điền)." 0
Set-Body12    "B5"  "- SKU: Mã số (Bắt buộc phải điền)." 0
Set-Body12    "C6"  "Điền giống như NAME." 1
Set-Body12    "B7"  "- PRICE: Giá cả (Có thể điền hoặc không điền)." 0

Set-Body12    "C12" "Điền" 2
Set-Body12    "C13" "VMS - MOBIPHONE." 2
Set-Body12    "C14" "VNM - VIETNAMOBILE." 2
Set-Body12    "C15" "VNP - VINAPHONE." 2
Set-Body12    "C16" "GSIM." 2
Set-Body12    "C17" "VTT - VIETTEL." 2

Set-Heading11 "B18" "- STOCK_QUANTITY: Số lượng (Có thể điền hoặc không điền)."
Set-Heading11 "B19" "- LEVEL: Hạng số. (Có thể điền hoặc không điền)."

Set-Body12    "C21" "NORMAL." 2
Set-Body12    "C22" "QUASI." 2
Set-Body12    "C23" "BRONZE." 2
Set-Body12    "C24" "SILVER." 2
Set-Body12    "C25" "GOLD." 2
Set-Body12    "C26" "PLATINUM" 2

Set-Heading12 "B2"  "- NAME: Số điện thoại/Số serial (Bắt buộc phải điền)."
Set-Heading12 "B8"  "- CATEGORY_ID: Loại sản phẩm (Bắt buộc điền)."

Set-Blank12   "B9"
Set-Body12    "C9"  "2 là serial SIM" 0

Set-BlankIndent "B10"
Set-Body12    "C10" "3 là số điện thoại" 0

Set-Heading12 "B11" "- BRAND: Tên nhà mạng (Bắt buộc điền)."

Set-BlankIndent "B12"
Set-Body12    "C20" "Điền" 1

# Row heights across the used range match the 12pt Times New Roman autofit.
$ws.Range("B2:C26").RowHeight = 15.75

# ---------------------------------------------------------------------
# 3. Products!I1 gets a new "LEVEL" header (last new shared string).
# ---------------------------------------------------------------------
$products.Range("I1").Value = "LEVEL"
$products.Range("J1").Select()

Write-Host "done"
